$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new paragraph right after the first paragraph (Heading1
#    "Play Dragon Spin free and enjoy the legendary wins") containing:
#      <empty run><bold "Meta description"><normal ": Read our review...">
# ------------------------------------------------------------------

$firstPara = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($firstPara.Range.End, $firstPara.Range.End)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dragon Spin online slot game and play for free. Enjoy oriental graphics and bonus rounds!</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($metaXml)

# InsertXML leaves behind one extra empty placeholder paragraph (needed to
# force a genuine paragraph break); remove it again.
$extraPara = $d.Paragraphs.Item(3)
$extraPara.Range.Delete() | Out-Null

# ------------------------------------------------------------------
# 2) Near the end of the document: remove the paragraph that duplicated
#    the bold "Play Dragon Spin free and enjoy the legendary wins" text,
#    and replace the text of the following italic paragraph with the new
#    "Prompt: ..." text (formatting / leading empty run untouched).
# ------------------------------------------------------------------

$boldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text -eq "Play Dragon Spin free and enjoy the legendary wins`r" -and $p.Range.Bold -eq -1) {
    $boldPara = $p
  }
}
$boldPara.Range.Delete() | Out-Null

$italicPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Italic -eq -1) {
    $italicPara = $p
  }
}

$iStart = $italicPara.Range.Start
$iEnd = $italicPara.Range.End
$italicTextRange = $d.Range($iStart, $iEnd - 1)
$italicTextRange.Text = 'Prompt: Design a cartoon-style image for the game "Dragon Spin" featuring a happy Maya warrior with glasses. The image should be eye-catching and vibrant, showcasing the mythical dragon theme of the game while also highlighting the fun and playful nature of the Maya warrior character. The warrior should be depicted with a large smile on their face, holding a staff or a sword and standing in a powerful pose. The background of the image should feature a cityscape inspired by ancient Mayan architecture, with a dragon flying in the distance. Color scheme should be vibrant and bold, incorporating shades of red, blue, and yellow. Overall, the image should capture the spirit of adventure and excitement that players can expect when playing "Dragon Spin."'
